$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Desicion_Tree)
$ws.Range("C2").Value = 94.84999999999999
$ws.Range("D2").Value = 98.06
$ws.Range("E2").Value = 97.87
$ws.Range("G2").Value = 96.5

# Row 4 (Deep_Neural_Network)
$ws.Range("C4").Value = 77.45
$ws.Range("D4").Value = 84.69
$ws.Range("E4").Value = 84.04000000000001
$ws.Range("F4").Value = 78.3
$ws.Range("G4").Value = 81
